$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the existing "spring" row (row 53), shifting
# all subsequent rows down by 9 (old row 53 -> new row 62, old row 110 -> new row 119).
$ws.Rows("53:61").Insert()

$newRows = @(
    @("Happy New Year", "あけましておめでとうございます"),
    @("Happy New Year (kanji)", "謹賀新年|きんがしんねん"),
    @("Thank you for all your kind help during the past year.", "昨年は大変お世話になりました|さくねんはたいへんおせわになりました"),
    @("I hope for your continued good will this year.", "本年もどうぞよろしくお願いいたします|ほんねんもどうぞよろしくおねがいいたします"),
    @("I hope you are keeping well during the hot weather.", "暑中お見舞い申し上げます|しょちゅうおみまいもうしあげます"),
    @("Congratulations on your graduation.", "ご卒業おめでとうございます|ごそつぎょうおめでとうございます"),
    @("Congratulations on your marriage.", "ご結婚おめでとうございます|ごけっこんおめでとうございます"),
    @("Happy Birthday", "誕生日おめでとう|たんじょうびおめでとう"),
    @("Get well soon.", "早くよくなってください|はやくよくなってください")
)

$r = 53
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
